$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.925285679134525
$ws.Range("D2").Value = 3.103869232618028
$ws.Range("E2").Value = 40.4747681248112
$ws.Range("F2").Value = 17.78066332051406
$ws.Range("G2").Value = 3.573216040938657
$ws.Range("M2").Value = 57.92686749379254
$ws.Range("O2").Value = 15.56101615178891

$ws.Range("C3").Value = 4.750027999926098
$ws.Range("D3").Value = 3.128693051900153
$ws.Range("E3").Value = 37.75599350956422
$ws.Range("F3").Value = 17.9591334604887
$ws.Range("G3").Value = 3.578045170979934
$ws.Range("M3").Value = 54.37435404903417
$ws.Range("O3").Value = 15.83860584521647

$ws.Range("C4").Value = 4.640900711807915
$ws.Range("D4").Value = 3.145869948088174
$ws.Range("E4").Value = 35.98604020278441
$ws.Range("F4").Value = 18.09174145983836
$ws.Range("G4").Value = 3.581124398193673
$ws.Range("M4").Value = 52.06825419896585
$ws.Range("O4").Value = 16.02637676403377

$ws.Range("C5").Value = 4.596125374968308
$ws.Range("D5").Value = 3.153342634848839
$ws.Range("E5").Value = 35.23929913919438
$ws.Range("F5").Value = 18.15133548150709
$ws.Range("G5").Value = 3.582408195568177
$ws.Range("M5").Value = 51.09715533469215
$ws.Range("O5").Value = 16.10708064282368

$ws.Range("C6").Value = 4.588674403562389
$ws.Range("D6").Value = 3.15461165017145
$ws.Range("E6").Value = 35.11376007494346
$ws.Range("F6").Value = 18.16155969861181
$ws.Range("G6").Value = 3.582623128178475
$ws.Range("M6").Value = 50.93401453297691
$ws.Range("O6").Value = 16.12072942357807

$ws.Range("C7").Value = 4.640297985192409
$ws.Range("D7").Value = 3.145968829805405
$ws.Range("E7").Value = 35.97607274510581
$ws.Range("F7").Value = 18.09252298244031
$ws.Range("G7").Value = 3.581141594188352
$ws.Range("M7").Value = 52.05528442758418
$ws.Range("O7").Value = 16.02744843713119

$ws.Range("C8").Value = 4.865216022931792
$ws.Range("D8").Value = 3.112018483005742
$ws.Range("E8").Value = 39.55811344960533
$ws.Range("F8").Value = 17.83727907595528
$ws.Range("G8").Value = 3.574857618052564
$ws.Range("M8").Value = 56.7278833872265
$ws.Range("O8").Value = 15.65302581684887

$ws.Range("C9").Value = 5.29119125856858
$ws.Range("D9").Value = 3.326766298068199
$ws.Range("E9").Value = 45.79200423058815
$ws.Range("F9").Value = 17.52992329139781
$ws.Range("G9").Value = 3.56342635894553
$ws.Range("M9").Value = 64.90090771461657
$ws.Range("O9").Value = 15.06449557629476

$ws.Range("C10").Value = 5.591353982844454
$ws.Range("D10").Value = 3.523301215962421
$ws.Range("E10").Value = 49.90124134556373
$ws.Range("F10").Value = 17.43593053576393
$ws.Range("G10").Value = 3.555550847333812
$ws.Range("M10").Value = 70.30617877219132
$ws.Range("O10").Value = 14.73294393748895

$ws.Range("C11").Value = 5.724495361806251
$ws.Range("D11").Value = 3.609433254402915
$ws.Range("E11").Value = 51.67087168169324
$ws.Range("F11").Value = 17.42466718375831
$ws.Range("G11").Value = 3.55207688074243
$ws.Range("M11").Value = 72.63655795396721
$ws.Range("O11").Value = 14.60681131080074

$ws.Range("C12").Value = 5.774378760615375
$ws.Range("D12").Value = 3.641579390379644
$ws.Range("E12").Value = 52.32685492261196
$ws.Range("F12").Value = 17.42515200557886
$ws.Range("G12").Value = 3.55077659641003
$ws.Range("M12").Value = 73.50068789788662
$ws.Range("O12").Value = 14.56284379958853

$ws.Range("C13").Value = 5.76365989656731
$ws.Range("D13").Value = 3.634677030576813
$ws.Range("E13").Value = 52.18620276792472
$ws.Range("F13").Value = 17.42483296705672
$ws.Range("G13").Value = 3.551055964944501
$ws.Range("M13").Value = 73.31539493152995
$ws.Range("O13").Value = 14.57214015146919

$ws.Range("C14").Value = 5.728610244655238
$ws.Range("D14").Value = 3.612087350813602
$ws.Range("E14").Value = 51.72512200579774
$ws.Range("F14").Value = 17.42461065794737
$ws.Range("G14").Value = 3.551969602625245
$ws.Range("M14").Value = 72.70801698065998
$ws.Range("O14").Value = 14.6031165585438

$ws.Range("C15").Value = 5.707070500047874
$ws.Range("D15").Value = 3.598189330272997
$ws.Range("E15").Value = 51.4408615526966
$ws.Range("F15").Value = 17.42509939663767
$ws.Range("G15").Value = 3.552531203688335
$ws.Range("M15").Value = 72.33359747748167
$ws.Range("O15").Value = 14.62259231475213

$ws.Range("C16").Value = 5.582580414738529
$ws.Range("D16").Value = 3.517606421803465
$ws.Range("E16").Value = 49.78360268917631
$ws.Range("F16").Value = 17.43731959360574
$ws.Range("G16").Value = 3.555780032426351
$ws.Range("M16").Value = 70.15130903484523
$ws.Range("O16").Value = 14.74170464502829

$ws.Range("C17").Value = 5.505304134390098
$ws.Range("D17").Value = 3.467332892204161
$ws.Range("E17").Value = 48.74153549956233
$ws.Range("F17").Value = 17.45304442719488
$ws.Range("G17").Value = 3.557800637902754
$ws.Range("M17").Value = 68.77972842144302
$ws.Range("O17").Value = 14.82127213698955

$ws.Range("C18").Value = 5.460538744903711
$ws.Range("D18").Value = 3.438108697845646
$ws.Range("E18").Value = 48.13277137941339
$ws.Range("F18").Value = 17.46504087750936
$ws.Range("G18").Value = 3.558973087395096
$ws.Range("M18").Value = 67.97872580048987
$ws.Range("O18").Value = 14.86934606065295

$ws.Range("C19").Value = 5.445328778116902
$ws.Range("D19").Value = 3.428160974487483
$ws.Range("E19").Value = 47.92503485646
$ws.Range("F19").Value = 17.4696027046975
$ws.Range("G19").Value = 3.559371830790851
$ws.Range("M19").Value = 67.7054368512236
$ws.Range("O19").Value = 14.88601227462633

$ws.Range("C20").Value = 5.513563593453192
$ws.Range("D20").Value = 3.472716494740209
$ws.Range("E20").Value = 48.85343622998003
$ws.Range("F20").Value = 17.45106333039227
$ws.Range("G20").Value = 3.557584482863358
$ws.Range("M20").Value = 68.92698743518507
$ws.Range("O20").Value = 14.8125613392852

$ws.Range("C21").Value = 5.738920008039332
$ws.Range("D21").Value = 3.618735235817497
$ws.Range("E21").Value = 51.86093469205524
$ws.Range("F21").Value = 17.42454535115195
$ws.Range("G21").Value = 3.551700834940255
$ws.Range("M21").Value = 72.88691492839068
$ws.Range("O21").Value = 14.59391300945234

$ws.Range("C22").Value = 5.883072017468903
$ws.Range("D22").Value = 3.711427742295451
$ws.Range("E22").Value = 53.74420054440342
$ws.Range("F22").Value = 17.43499268027773
$ws.Range("G22").Value = 3.547944114141152
$ws.Range("M22").Value = 75.36817331431817
$ws.Range("O22").Value = 14.47326048848868

$ws.Range("C23").Value = 5.806435249338155
$ws.Range("D23").Value = 3.662206027761803
$ws.Range("E23").Value = 52.74653093231685
$ws.Range("F23").Value = 17.42680585784758
$ws.Range("G23").Value = 3.549941176280787
$ws.Range("M23").Value = 74.05359637252666
$ws.Range("O23").Value = 14.53553611584481

$ws.Range("C24").Value = 5.509830542046719
$ws.Range("D24").Value = 3.470283568734412
$ws.Range("E24").Value = 48.80287603144036
$ws.Range("F24").Value = 17.45194979382599
$ws.Range("G24").Value = 3.557682172961408
$ws.Range("M24").Value = 68.8604504641953
$ws.Range("O24").Value = 14.81649225351308

$ws.Range("C25").Value = 5.177966517959183
$ws.Range("D25").Value = 3.251496742913989
$ws.Range("E25").Value = 44.18906665281732
$ws.Range("F25").Value = 17.59092136329062
$ws.Range("G25").Value = 3.566425423818949
$ws.Range("M25").Value = 62.79558542575236
$ws.Range("O25").Value = 14.8125613392852
